$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 74082.25
$ws.Range("J81").Value = 74082.25
$ws.Range("L81").Value = 74082.25
$ws.Range("N81").Value = -76078.25
$ws.Range("H84").Value = 74082.25
$ws.Range("J84").Value = 74082.25
$ws.Range("L84").Value = 222246.75
$ws.Range("N84").Value = -232230.75
$ws.Range("H86").Value = 6026.45
$ws.Range("J86").Value = 8095
$ws.Range("L86").Value = 8095
$ws.Range("N86").Value = -10341
$ws.Range("H89").Value = 6026.45
$ws.Range("J89").Value = 8095
$ws.Range("L89").Value = 40475
$ws.Range("N89").Value = -51707
$ws.Range("H106").Value = 3476.1538
$ws.Range("I106").Value = 2924.5
$ws.Range("J106").Value = 3721.3333
$ws.Range("K106").Value = 2924.5
$ws.Range("L106").Value = 3721.3333
$ws.Range("M106").Value = -2293.5
$ws.Range("N106").Value = -4983.3333
$ws.Range("H116").Value = 6219.4375
$ws.Range("I116").Value = 5513.778
$ws.Range("K116").Value = 5513.778
$ws.Range("M116").Value = -2071.778

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1184.5
$ws.Range("I2").Value = 1184.5
$ws.Range("K2").Value = 1184.5
$ws.Range("M2").Value = -1071.5
$ws.Range("H5").Value = 83.111115
$ws.Range("I5").Value = 84
$ws.Range("K5").Value = 84
$ws.Range("M5").Value = 28
$ws.Range("H31").Value = 43793.25
$ws.Range("I31").Value = 3799.6
$ws.Range("J31").Value = 110449.336
$ws.Range("K31").Value = 3799.6
$ws.Range("L31").Value = 110449.336
$ws.Range("M31").Value = -3505.6
$ws.Range("N31").Value = -111037.336
$ws.Range("H45").Value = 2282.8667
$ws.Range("I45").Value = 1848.6
$ws.Range("K45").Value = 1848.6
$ws.Range("M45").Value = -1471.6
$ws.Range("H74").Value = 5819874.5
$ws.Range("I74").Value = 6758382
$ws.Range("K74").Value = 6758382
$ws.Range("M74").Value = -6757508
$ws.Range("H77").Value = 5819874.5
$ws.Range("I77").Value = 6758382
$ws.Range("K77").Value = 33791910
$ws.Range("M77").Value = -33787542
$ws.Range("H116").Value = 1184.5
$ws.Range("I116").Value = 1184.5
$ws.Range("K116").Value = 1184.5
$ws.Range("M116").Value = 1109.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1184.5
$ws.Range("I3").Value = 1184.5
$ws.Range("K3").Value = 1184.5
$ws.Range("M3").Value = -1070.5
$ws.Range("H4").Value = 83.111115
$ws.Range("I4").Value = 84
$ws.Range("K4").Value = 84
$ws.Range("M4").Value = 31
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = $null
$ws.Range("H11").Value = 17166.666
$ws.Range("I11").Value = 750
$ws.Range("K11").Value = 750
$ws.Range("M11").Value = -610
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = $null
$ws.Range("H96").Value = 42670
$ws.Range("J96").Value = 82897
$ws.Range("L96").Value = 82897
$ws.Range("N96").Value = -88389
$ws.Range("H99").Value = 1941.7858
$ws.Range("I99").Value = 1448.75
$ws.Range("K99").Value = 1448.75
$ws.Range("M99").Value = 49.25
$ws.Range("H105").Value = 1577.125
$ws.Range("I105").Value = 1445.2858
$ws.Range("K105").Value = 1445.2858
$ws.Range("M105").Value = 301.7141999999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 51451.43
$ws.Range("I2").Value = 150
$ws.Range("J2").Value = 60001.668
$ws.Range("K2").Value = 150
$ws.Range("L2").Value = 60001.668
$ws.Range("M2").Value = -37
$ws.Range("N2").Value = -60227.668
$ws.Range("H3").Value = 1999.5
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1999.5
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 1999.5
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = -2225.5
$ws.Range("H58").Value = 2585.1667
$ws.Range("I58").Value = 2882.5
$ws.Range("K58").Value = 2882.5
$ws.Range("M58").Value = -2679.5
$ws.Range("H94").Value = 4351.8335
$ws.Range("I94").Value = 350
$ws.Range("J94").Value = 4715.636
$ws.Range("K94").Value = 350
$ws.Range("L94").Value = 4715.636
$ws.Range("M94").Value = 101
$ws.Range("N94").Value = -5617.636
$ws.Range("H136").Value = 2585.1667
$ws.Range("I136").Value = 2882.5
$ws.Range("K136").Value = 8647.5
$ws.Range("M136").Value = -6097.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16350073
$ws.Range("I4").Value = 22400026
$ws.Range("J4").Value = 10300120
$ws.Range("K4").Value = 67200078
$ws.Range("L4").Value = 30900360
$ws.Range("M4").Value = -67199966
$ws.Range("N4").Value = -30900584
$ws.Range("H14").Value = 567.5
$ws.Range("I14").Value = 567.5
$ws.Range("K14").Value = 1702.5
$ws.Range("M14").Value = -1529.5
$ws.Range("H34").Value = 2410.3572
$ws.Range("J34").Value = 4285.5713
$ws.Range("L34").Value = 12856.7139
$ws.Range("N34").Value = -13024.7139
$ws.Range("H46").Value = 440.55554
$ws.Range("J46").Value = 527.2
$ws.Range("L46").Value = 1581.6
$ws.Range("N46").Value = -1763.6
$ws.Range("H74").Value = 9350
$ws.Range("J74").Value = 15000
$ws.Range("L74").Value = 45000
$ws.Range("N74").Value = -47122
$ws.Range("H77").Value = 9350
$ws.Range("J77").Value = 15000
$ws.Range("L77").Value = 135000
$ws.Range("N77").Value = -145608
$ws.Range("H80").Value = 3999
$ws.Range("J80").Value = 3999
$ws.Range("L80").Value = 11997
$ws.Range("N80").Value = -13869
$ws.Range("H81").Value = 1756.5
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("H83").Value = 3999
$ws.Range("J83").Value = 3999
$ws.Range("L83").Value = 35991
$ws.Range("N83").Value = -45351
$ws.Range("H84").Value = 1756.5
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").Value = $null
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").Value = $null

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2578.9167
$ws.Range("J102").Value = 2250
$ws.Range("L102").Value = 2250
$ws.Range("N102").Value = -5494

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3665.611
$ws.Range("J40").Value = 4909.273
$ws.Range("L40").Value = 4909.273
$ws.Range("N40").Value = -5181.273
